# MitsosBarton2006Ex323 nonconvex experiment generator - update M_Stationary
# (alpha_zero) results: x/y moved to a new stationary point (x=2.09, y=3.8)
# and every dependent expression-evaluation / multiplier / vector value is
# recomputed accordingly. Sheet "Funciones_Objetivo" (expressions, x & y
# symbol headers) and "Vector_Alpha" stay untouched, matching the commit.
#
# All the written values are numeric-looking text (the generator always
# stored numbers as text/shared-strings, never as native numbers), so each
# target range is switched to Text format before the value is typed in,
# then ClearFormats() drops the now-superfluous "@" format again so the
# cells end up with no cell style at all - exactly like the source file.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Sheet, $RangeAddress, $TextValue)
    $rng = $Sheet.Range($RangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $TextValue
    $rng.ClearFormats()
}

# NOTE: sheet indices are used (not names) because the workbook has two
# sheets whose names differ only by case - "Vector_bf" and "Vector_BF" -
# and name-based Item() lookup is case-insensitive, so it would resolve
# both to the same (wrong) sheet.

# --- Restricciones_del_lider (sheet index 2) ------------------------------
$ws2 = $wb.Worksheets.Item(2)

Set-TextValue $ws2 "A2" "2.09 - x"
Set-TextValue $ws2 "B2" "-3.09"
Set-TextValue $ws2 "D2" "0.86"

Set-TextValue $ws2 "A3" "-2.09 + x"
Set-TextValue $ws2 "B3" "1.0899999999999999"
Set-TextValue $ws2 "D3" "0.62"

Set-TextValue $ws2 "A4" "41.02289999999999 + x - y - 9(x^2)"
Set-TextValue $ws2 "B4" "-40.02289999999999"
Set-TextValue $ws2 "D4" "0.58"

# --- Restricciones_del_follower (sheet index 3) ---------------------------
$ws3 = $wb.Worksheets.Item(3)

Set-TextValue $ws3 "A2" "-22.9596 + (-0.5 + x)*(y^2)"
Set-TextValue $ws3 "B2" "22.9596"
Set-TextValue $ws3 "D2" "0.69"
Set-TextValue $ws3 "E2" "6.4"
Set-TextValue $ws3 "F2" "8.100000000000001"

Set-TextValue $ws3 "A3" "-3.8 + y"
Set-TextValue $ws3 "B3" "2.8"
Set-TextValue $ws3 "D3" "0.65"
Set-TextValue $ws3 "E3" "4.4"
Set-TextValue $ws3 "F3" "3.5"

Set-TextValue $ws3 "A4" "-5.8 - y"
Set-TextValue $ws3 "B4" "-4.8"
Set-TextValue $ws3 "D4" "0.32"
Set-TextValue $ws3 "E4" "5.2"
Set-TextValue $ws3 "F4" "6.2"

# --- Punto_modificado (sheet index 4) --------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Set-TextValue $ws4 "A2" "2.09"
Set-TextValue $ws4 "B2" "3.8"

# --- Vector_bf (sheet index 5) ----------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-TextValue $ws5 "A2" "-9.667959999999997"

# --- Vector_BF (sheet index 6) -----------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-TextValue $ws6 "A2" "-71.9364"
Set-TextValue $ws6 "A3" "-75.9576"

# Vector_Alpha (sheet index 7) and Funciones_Objetivo (sheet index 1) are
# unchanged.
